$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - CON values (B2:E2)
$ws.Range("B2").Value = 47.411870575783915
$ws.Range("C2").Value = 54.726252459006133
$ws.Range("D2").Value = 50.707563859743551
$ws.Range("E2").Value = 55.149544983943933

# Row 3 - STR values (B3:E3)
$ws.Range("B3").Value = 44.659973050356776
$ws.Range("C3").Value = 45.107778304367976
$ws.Range("D3").Value = 44.245052626267544
$ws.Range("E3").Value = 53.747993844062918

# Update the selection to match the new reduced range
$ws.Range("B1:E3").Select()
